$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 171
$ws.Range("J2").Value = 807
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 215
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = 157
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 89
$ws.Range("T2").Value = 129
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 1172
$ws.Range("X2").Value = 1141
$ws.Range("Z2").Value = 14
$ws.Range("AA2").Value = 7
